$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Footer "default" (index 1) -> footer2.xml: Pearson logo, wp:docPr id="2" ---
$ftrDefault = $sec.Footers.Item(1)
$rngFtrDefault = $ftrDefault.Range
$xmlFtrDefault = $rngFtrDefault.WordOpenXML
$xmlFtrDefault = $xmlFtrDefault.Replace('descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image2.png"', 'descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image1.png"')
$xmlFtrDefault = $xmlFtrDefault.Replace('descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image2.png"', 'descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"')
$rngFtrDefault.WordOpenXML = $xmlFtrDefault

# --- Footer "first page" (index 2) -> footer1.xml: Pearson logo, wp:docPr id="3" ---
$ftrFirst = $sec.Footers.Item(2)
$rngFtrFirst = $ftrFirst.Range
$xmlFtrFirst = $rngFtrFirst.WordOpenXML
$xmlFtrFirst = $xmlFtrFirst.Replace('descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="3" name="image2.png"', 'descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="3" name="image1.png"')
$xmlFtrFirst = $xmlFtrFirst.Replace('descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image2.png"', 'descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"')
$rngFtrFirst.WordOpenXML = $xmlFtrFirst

# --- Header "first page" (index 2) -> header1.xml: BTec logo, wp:docPr id="1" ---
$hdrFirst = $sec.Headers.Item(2)
$rngHdrFirst = $hdrFirst.Range
$xmlHdrFirst = $rngHdrFirst.WordOpenXML
$xmlHdrFirst = $xmlHdrFirst.Replace('descr="BTec_Logo-Orange" id="1" name="image1.jpg"', 'descr="BTec_Logo-Orange" id="1" name="image2.jpg"')
$xmlHdrFirst = $xmlHdrFirst.Replace('descr="BTec_Logo-Orange" id="0" name="image1.jpg"', 'descr="BTec_Logo-Orange" id="0" name="image2.jpg"')
$rngHdrFirst.WordOpenXML = $xmlHdrFirst

Write-Host "Renamed inline shapes in footers/header"
